$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2025-12-25 01:02:23","Admin","Login","login_success","Role: admin"),
    @("2025-12-25 01:02:23","Admin","dashboard","access_granted","Opened dashboard page"),
    @("2025-12-25 01:02:25","Admin","quotation","access_granted","Opened quotation page"),
    @("2025-12-25 01:02:38","Admin","invoice","access_granted","Opened invoice page"),
    @("2025-12-25 01:06:51","Admin","Login","login_success","Role: admin"),
    @("2025-12-25 01:06:51","Admin","dashboard","access_granted","Opened dashboard page"),
    @("2025-12-25 01:06:53","Admin","quotation","access_granted","Opened quotation page"),
    @("2025-12-25 01:06:55","Admin","invoice","access_granted","Opened invoice page"),
    @("2025-12-25 01:07:00","Admin","invoice","access_granted","Opened invoice page"),
    @("2025-12-25 01:07:00","Admin","invoice","access_granted","Opened invoice page"),
    @("2025-12-25 01:07:03","Admin","invoice","access_granted","Opened invoice page"),
    @("2025-12-25 01:07:07","Admin","invoice","access_granted","Opened invoice page"),
    @("2025-12-25 01:07:07","Admin","invoice","access_granted","Opened invoice page"),
    @("2025-12-25 01:07:11","Admin","invoice","access_granted","Opened invoice page"),
    @("2025-12-25 01:07:19","Admin","invoice","access_granted","Opened invoice page"),
    @("2025-12-25 01:07:22","Admin","invoice","access_granted","Opened invoice page"),
    @("2025-12-25 01:07:24","Admin","invoice","access_granted","Opened invoice page")
)

$startRow = 27
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
    $ws.Cells.Item($row, 5).Value = $data[$i][4]
}
